# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (string) storage for the Price column cells we are about to
# rewrite, so numeric-looking strings (e.g. "295.18", "1.00") are kept as
# literal text instead of being auto-coerced to numbers (matches the
# original inlineStr cell type). NumberFormat must be set per-cell since
# this COM bridge only applies Range.NumberFormat to the first area of a
# multi-area/union range.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values exactly as they appear in the updated sheet.
$ws.Range('D2').Value = '39.374.63'
$ws.Range('E2').Value = '  -1.85%  '
$ws.Range('D3').Value = '2.197.12'
$ws.Range('E3').Value = '  -6.07%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '295.18'
$ws.Range('E5').Value = '  -4.27%  '
$ws.Range('D6').Value = '81.69'
$ws.Range('E6').Value = '  -4.06%  '
$ws.Range('D7').Value = '0.511'
$ws.Range('E7').Value = '  -3.63%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').Value = '0.467'
$ws.Range('E9').Value = '  -3.91%  '
$ws.Range('D10').Value = '0.0769'
$ws.Range('E10').Value = '  -6.48%  '
$ws.Range('D11').Value = '29.12'
$ws.Range('E11').Value = '  -3.96%  '
$ws.Range('E12').Value = '  -10.68%  '
$ws.Range('D14').Value = '2.536.54'
$ws.Range('E14').Value = '  -5.91%  '
$ws.Range('D15').Value = '6.23'
$ws.Range('E15').Value = '  -3.35%  '
$ws.Range('D16').Value = '13.95'
$ws.Range('E16').Value = '  -5.64%  '
$ws.Range('D17').Value = '2.201.48'
$ws.Range('E17').Value = '  -5.51%  '
$ws.Range('D18').Value = '0.711'
$ws.Range('E18').Value = '  -5.44%  '
$ws.Range('D19').Value = '39.268.62'
$ws.Range('E19').Value = '  -1.86%  '
$ws.Range('D20').Value = '0.0₃0869'
$ws.Range('E20').Value = '  -4.20%  '
$ws.Range('D21').Value = '5.71'
$ws.Range('E21').Value = '  -6.35%  '
$ws.Range('D22').Value = '64.75'
$ws.Range('E22').Value = '  -4.51%  '
$ws.Range('D23').Value = '10.28'
$ws.Range('E23').Value = '  -4.18%  '
$ws.Range('D24').Value = '225.15'
$ws.Range('E24').Value = '  -4.50%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '2.40'
$ws.Range('E26').Value = '  -6.05%  '
$ws.Range('D27').Value = '1.80'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').Value = '22.55'
$ws.Range('E28').Value = '  -4.21%  '
$ws.Range('E29').Value = '  +0.63%  '
$ws.Range('D30').Value = '9.07'
$ws.Range('E30').Value = '  -1.69%  '
$ws.Range('D31').Value = '148.56'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('D32').Value = '31.65'
$ws.Range('E32').Value = '  -9.89%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  -6.74%  '
$ws.Range('D35').Value = '0.0694'
$ws.Range('E35').Value = '  -4.50%  '
$ws.Range('E36').Value = '  -4.81%  '
$ws.Range('D37').Value = '0.110'
$ws.Range('E37').Value = '  -3.56%  '
$ws.Range('D38').Value = '15.38'
$ws.Range('E38').Value = '  -3.01%  '
$ws.Range('D39').Value = '0.0958'
$ws.Range('E39').Value = '  -4.38%  '
$ws.Range('D40').Value = '2.63'
$ws.Range('E40').Value = '  -5.48%  '
$ws.Range('D41').Value = '1.64'
$ws.Range('E41').Value = '  -4.30%  '
$ws.Range('D42').Value = '3.60'
$ws.Range('E42').Value = '  -5.72%  '
$ws.Range('D43').Value = '1.896.32'
$ws.Range('E43').Value = '  -2.63%  '
$ws.Range('E44').Value = '  -9.89%  '
$ws.Range('D45').Value = '0.0259'
$ws.Range('E45').Value = '  -2.98%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '16.05'
$ws.Range('E46').Value = '  -9.94%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '8.93'
$ws.Range('E47').Value = '  -4.02%  '
$ws.Range('D48').Value = '2.59'
$ws.Range('E48').Value = '  -3.49%  '
$ws.Range('D49').Value = '71.31'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').Value = '2.403.89'
$ws.Range('E50').Value = '  -5.79%  '
$ws.Range('D51').Value = '87.05'
$ws.Range('E51').Value = '  -6.30%  '

Write-Host "Applied cryptos.xlsx update"
